# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) values for the
# 9c330023-... row (row 3) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-26 06:35:27"
$wsZhCn.Range("G3").Value = "2016-02-26 06:36:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-26 06:35:42"
$wsDeDe.Range("G3").Value = "2016-02-26 06:36:47"
